$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.312.51'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '3.568.52'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  -0.05%  '
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '605.20'
$cell.Style = $origStyle
$ws.Range("E5").Value = '  -0.14%  '
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '147.46'
$cell.Style = $origStyle
$ws.Range("E6").Value = '  +2.11%  '
$ws.Range("D7").Value = '3.568.32'
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("E8").Value = '  -0.11%  '
$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.490'
$cell.Style = $origStyle
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("E10").Value = '  -1.33%  '
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.91'
$cell.Style = $origStyle
$ws.Range("E11").Value = '  +1.41%  '
$ws.Range("E12").Value = '  -0.73%  '
$ws.Range("D13").Value = '4.175.25'
$ws.Range("E13").Value = '  +0.15%  '
$ws.Range("E14").Value = '  -1.45%  '
$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '29.39'
$cell.Style = $origStyle
$ws.Range("E15").Value = '  -3.38%  '
$ws.Range("D16").Value = '3.572.80'
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("E17").Value = '  +1.62%  '
$ws.Range("D18").Value = '66.291.88'
$ws.Range("E18").Value = '  +0.01%  '
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.97'
$cell.Style = $origStyle
$ws.Range("E19").Value = '  -3.60%  '
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.29'
$cell.Style = $origStyle
$ws.Range("E20").Value = '  +1.09%  '
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '14.74'
$cell.Style = $origStyle
$ws.Range("E21").Value = '  -0.44%  '
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '421.02'
$cell.Style = $origStyle
$ws.Range("E22").Value = '  -2.25%  '
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.609'
$cell.Style = $origStyle
$ws.Range("E23").Value = '  -0.87%  '
$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '77.85'
$cell.Style = $origStyle
$ws.Range("E24").Value = '  -2.05%  '
$ws.Range("D25").Value = '3.708.33'
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  +0.65%  '
$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '9.30'
$cell.Style = $origStyle
$ws.Range("E28").Value = '  +1.42%  '
$ws.Range("E29").Value = '  +1.18%  '
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.49'
$cell.Style = $origStyle
$ws.Range("E30").Value = '  -0.52%  '
$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = $origStyle
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("D32").Value = '3.565.13'
$ws.Range("E32").Value = '  +0.17%  '
$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.157'
$cell.Style = $origStyle
$ws.Range("E33").Value = '  +4.06%  '
$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '24.80'
$cell.Style = $origStyle
$ws.Range("E34").Value = '  -2.52%  '
$ws.Range("E35").Value = '  -3.31%  '
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.71'
$cell.Style = $origStyle
$ws.Range("E37").Value = '  -1.83%  '
$ws.Range("E38").Value = '  -3.24%  '
$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.63'
$cell.Style = $origStyle
$ws.Range("E39").Value = '  -4.86%  '
$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '175.24'
$cell.Style = $origStyle
$ws.Range("E40").Value = '  +0.37%  '
$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0842'
$cell.Style = $origStyle
$ws.Range("E41").Value = '  -1.07%  '
$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.16'
$cell.Style = $origStyle
$ws.Range("E42").Value = '  -0.76%  '
$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.874'
$cell.Style = $origStyle
$ws.Range("E43").Value = '  -1.49%  '
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '45.82'
$cell.Style = $origStyle
$ws.Range("E44").Value = '  -0.27%  '
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.84'
$cell.Style = $origStyle
$ws.Range("E45").Value = '  -4.44%  '
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '23.47'
$cell.Style = $origStyle
$ws.Range("E48").Value = '  +0.29%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '24.16'
$cell.Style = $origStyle
$ws.Range("E49").Value = '  -3.00%  '
$ws.Range("B50").Value = 'ONDO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.13'
$cell.Style = $origStyle
$ws.Range("E50").Value = '  -5.18%  '
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.10'
$cell.Style = $origStyle
$ws.Range("E51").Value = '  -0.55%  '
